$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 'Environmental Impacts Assessment'
$ws.Range("C4").Value = 'Environmental Impacts Assessment'
$ws.Range("B6").Value = '''1'
$ws.Range("C6").Value = '''1'
$ws.Range("B7").Value = '90 h'
$ws.Range("C7").Value = '90 h'
$ws.Range("B8").Value = '''01/01/2022'
$ws.Range("C8").Value = '''01/01/2022'
$ws.Range("B10").Value = 'Propiciar aos alunos conhecimento sobre os fundamentos, objetivos e métodos da Avaliação de Impacto Ambiental.'
$ws.Range("C10").Value = 'Propiciar aos alunos conhecimento sobre os fundamentos, objetivos e métodos da Avaliação de Impacto Ambiental.'
$ws.Range("B11").Value = 'Provide knowledge on the fundamentals, objectives and methods of environmental impact assessment.'
$ws.Range("C11").Value = 'Provide knowledge on the fundamentals, objectives and methods of environmental impact assessment.'
$ws.Range("B14").Value = 'Os objetivos da Avaliação de Impacto Ambiental (AIA). O processo da AIA: triagem, definição do escopo, estudos de base, análise de impactos ambientais, mitigação, análise técnica e acompanhamento. Requisitos legais.'
$ws.Range("C14").Value = 'Os objetivos da Avaliação de Impacto Ambiental (AIA). O processo da AIA: triagem, definição do escopo, estudos de base, análise de impactos ambientais, mitigação, análise técnica e acompanhamento. Requisitos legais.'
$ws.Range("B16").Value = 'Conceitos básicos e definições. Origem e difusão da AIA. AIA e licenciamento: objetivos e fundamentos. Quadro legal e institucional brasileiro para a AIA. O processo de AIA e seus componentes. Etapas do planejamento e execução de um Estudo de Impacto Ambiental. Alternativas tecnológicas e de localização. Estudos de base e diagnóstico ambiental. Técnicas de identificação e previsão de impactos. Métodos e critérios de avaliação da importância dos impactos. Plano de gestão ambiental: medidas mitigadoras, compensatórias, de valorização e monitoramento. Tomada de decisão e acompanhamento. Estudos de caso.'
$ws.Range("C16").Value = 'Conceitos básicos e definições. Origem e difusão da AIA. AIA e licenciamento: objetivos e fundamentos. Quadro legal e institucional brasileiro para a AIA. O processo de AIA e seus componentes. Etapas do planejamento e execução de um Estudo de Impacto Ambiental. Alternativas tecnológicas e de localização. Estudos de base e diagnóstico ambiental. Técnicas de identificação e previsão de impactos. Métodos e critérios de avaliação da importância dos impactos. Plano de gestão ambiental: medidas mitigadoras, compensatórias, de valorização e monitoramento. Tomada de decisão e acompanhamento. Estudos de caso.'
$ws.Range("B17").Value = 'Basic concepts and definitions. Origin and dissemination of the Environmental Impact Assessment (EIA). EIA and licensing: objectives and fundamentals. Brazilian legal and institutional framework for EIA. The EIA process and its components. Steps in planning and execution of an Environmental Impact Statement. Technological and localization alternatives. Baseline studies. Impact identification and prediction techniques. Methods and criteria for determining impact significance. Environmental management plan: mitigation, compensation, monitoring. Decision making and follow-up. Case studies.'
$ws.Range("C17").Value = 'Basic concepts and definitions. Origin and dissemination of the Environmental Impact Assessment (EIA). EIA and licensing: objectives and fundamentals. Brazilian legal and institutional framework for EIA. The EIA process and its components. Steps in planning and execution of an Environmental Impact Statement. Technological and localization alternatives. Baseline studies. Impact identification and prediction techniques. Methods and criteria for determining impact significance. Environmental management plan: mitigation, compensation, monitoring. Decision making and follow-up. Case studies.'
$ws.Range("B19").Value = 'Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas.'
$ws.Range("C19").Value = 'Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas.'
$ws.Range("B20").Value = 'Média ponderada de atividades e provas.'
$ws.Range("C20").Value = 'Média ponderada de atividades e provas.'
$ws.Range("B21").Value = '1 (uma) prova escrita'
$ws.Range("C21").Value = '1 (uma) prova escrita'
$ws.Range("B22").Value = 'Bibliografia básicaSÁNCHEZ, L.E., Avaliação de impacto ambiental: conceitos e métodos. Oficina de textos: São Paulo, 2013. 583p.CALIJURI, M.C., CUNHA, D.G.F. (Org.), Engenharia ambiental: conceitos, tecnologia e gestão. Elsevier: Rio de Janeiro, 2019. 685p.Bibliografia complementar:COMPANHIA AMBIENTAL DO ESTADO DE SÃO PAULO - CETESB. MANUAL PARA ELABORAÇÃO DE ESTUDOS PARA O LICENCIAMENTO COM AVALIAÇÃO DE IMPACTO AMBIENTAL. Departamento de Desenvolvimento de Ações Estratégicas para o Licenciamento da Diretoria I- ID - CETESB. Anexo único, 2014. 250p.'
$ws.Range("C22").Value = 'Bibliografia básicaSÁNCHEZ, L.E., Avaliação de impacto ambiental: conceitos e métodos. Oficina de textos: São Paulo, 2013. 583p.CALIJURI, M.C., CUNHA, D.G.F. (Org.), Engenharia ambiental: conceitos, tecnologia e gestão. Elsevier: Rio de Janeiro, 2019. 685p.Bibliografia complementar:COMPANHIA AMBIENTAL DO ESTADO DE SÃO PAULO - CETESB. MANUAL PARA ELABORAÇÃO DE ESTUDOS PARA O LICENCIAMENTO COM AVALIAÇÃO DE IMPACTO AMBIENTAL. Departamento de Desenvolvimento de Ações Estratégicas para o Licenciamento da Diretoria I- ID - CETESB. Anexo único, 2014. 250p.'
